$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("E3").Value = 21
$ws.Range("F3").Value = 8
$ws.Range("H3").Value = 8

# Row 8
$ws.Range("E8").Value = 38

# Row 10
$ws.Range("E10").Value = 21
$ws.Range("F10").Value = 9
$ws.Range("H10").Value = 9

# Row 15
$ws.Range("E15").Value = 83

# Row 16
$ws.Range("E16").Value = 284

# Row 17
$ws.Range("E17").Value = 18
$ws.Range("F17").Value = 11
$ws.Range("H17").Value = 11

$wb.Save()
